$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 256.21738
$ws.Range("I39").Value = 41.23077
$ws.Range("J39").Value = 535.7
$ws.Range("K39").Value = 123.69231
$ws.Range("L39").Value = 1607.1
$ws.Range("M39").Value = 172.30769
$ws.Range("N39").Value = -2199.1
$ws.Range("H74").Value = 7989.8887
$ws.Range("I74").Value = 9800.799999999999
$ws.Range("K74").Value = 9800.799999999999
$ws.Range("M74").Value = -8864.799999999999
$ws.Range("H77").Value = 7989.8887
$ws.Range("I77").Value = 9800.799999999999
$ws.Range("K77").Value = 49004
$ws.Range("M77").Value = -44324
$ws.Range("H80").Value = 16821.125
$ws.Range("I80").Value = 1100
$ws.Range("J80").Value = 22061.5
$ws.Range("K80").Value = 3300
$ws.Range("L80").Value = 66184.5
$ws.Range("M80").Value = -2302
$ws.Range("N80").Value = -68180.5
$ws.Range("H83").Value = 16821.125
$ws.Range("I83").Value = 1100
$ws.Range("J83").Value = 22061.5
$ws.Range("K83").Value = 9900
$ws.Range("L83").Value = 198553.5
$ws.Range("M83").Value = -4908
$ws.Range("N83").Value = -208537.5
$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 898
$ws.Range("N99").ClearContents()
$ws.Range("H111").Value = 555.5
$ws.Range("I111").Value = 470.2
$ws.Range("J111").Value = 768.75
$ws.Range("K111").Value = 1410.6
$ws.Range("L111").Value = 2306.25
$ws.Range("M111").Value = 1656.4
$ws.Range("N111").Value = -8440.25
$ws.Range("H133").Value = 73968.42999999999
$ws.Range("J133").Value = 73968.42999999999
$ws.Range("L133").Value = 73968.42999999999
$ws.Range("N133").Value = -84088.42999999999
$ws.Range("H134").Value = 98971.42999999999
$ws.Range("J134").Value = 98971.42999999999
$ws.Range("L134").Value = 98971.42999999999
$ws.Range("N134").Value = -109111.43
$ws.Range("H136").Value = 77977.14
$ws.Range("J136").Value = 77977.14
$ws.Range("L136").Value = 77977.14
$ws.Range("N136").Value = -88177.14
$ws.Range("H137").Value = 634685.1
$ws.Range("I137").Value = 4085.4167
$ws.Range("K137").Value = 12256.2501
$ws.Range("M137").Value = -9706.250100000001
$ws.Range("H138").Value = 2764.818
$ws.Range("I138").Value = 2210.5
$ws.Range("J138").Value = 2888
$ws.Range("K138").Value = 6631.5
$ws.Range("L138").Value = 8664
$ws.Range("M138").Value = -1491.5
$ws.Range("N138").Value = -18944

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 65825
$ws.Range("I74").Value = 127164.75
$ws.Range("K74").Value = 127164.75
$ws.Range("M74").Value = -126290.75
$ws.Range("H77").Value = 65825
$ws.Range("I77").Value = 127164.75
$ws.Range("K77").Value = 635823.75
$ws.Range("M77").Value = -631455.75
$ws.Range("H94").Value = 34999.5
$ws.Range("J94").Value = 34999.5
$ws.Range("L94").Value = 34999.5
$ws.Range("N94").Value = -36801.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 94995.375
$ws.Range("J108").Value = 94995.375
$ws.Range("L108").Value = 94995.375
$ws.Range("N108").Value = -102675.375
$ws.Range("H119").Value = 98771.125
$ws.Range("J119").Value = 98771.125
$ws.Range("L119").Value = 98771.125
$ws.Range("N119").Value = -108447.125
$ws.Range("H135").Value = 98983.28999999999
$ws.Range("J135").Value = 98983.28999999999
$ws.Range("L135").Value = 98983.28999999999
$ws.Range("N135").Value = -109123.29
$ws.Range("H138").Value = 76664.44500000001
$ws.Range("J138").Value = 76664.44500000001
$ws.Range("L138").Value = 76664.44500000001
$ws.Range("N138").Value = -86944.44500000001
$ws.Range("H140").Value = 43498.934
$ws.Range("J140").Value = 43498.934
$ws.Range("L140").Value = 43498.934
$ws.Range("N140").Value = -53858.934

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3883.3157
$ws.Range("I31").Value = 2403.75
$ws.Range("J31").Value = 4959.364
$ws.Range("K31").Value = 2403.75
$ws.Range("L31").Value = 4959.364
$ws.Range("M31").Value = -2108.75
$ws.Range("N31").Value = -5549.364
$ws.Range("H34").Value = 3883.3157
$ws.Range("I34").Value = 2403.75
$ws.Range("J34").Value = 4959.364
$ws.Range("K34").Value = 2403.75
$ws.Range("L34").Value = 4959.364
$ws.Range("M34").Value = -2201.75
$ws.Range("N34").Value = -5363.364
$ws.Range("H114").Value = 63746
$ws.Range("J114").Value = 63746
$ws.Range("L114").Value = 63746
$ws.Range("N114").Value = -72424
$ws.Range("H138").Value = 54353.332
$ws.Range("J138").Value = 54897.5
$ws.Range("L138").Value = 54897.5
$ws.Range("N138").Value = -65177.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 16325.8
$ws.Range("I7").Value = 20219.75
$ws.Range("J7").Value = 750
$ws.Range("K7").Value = 60659.25
$ws.Range("L7").Value = 2250
$ws.Range("M7").Value = -60547.25
$ws.Range("N7").Value = -2474
$ws.Range("H63").Value = 300
$ws.Range("I63").Value = 300
$ws.Range("K63").Value = 900
$ws.Range("M63").Value = -151
$ws.Range("H66").Value = 300
$ws.Range("I66").Value = 300
$ws.Range("K66").Value = 2700
$ws.Range("M66").Value = 1044
$ws.Range("H68").Value = 85223.586
$ws.Range("J68").Value = 92743.91
$ws.Range("L68").Value = 278231.73
$ws.Range("N68").Value = -279853.73
$ws.Range("H69").Value = 9750
$ws.Range("J69").Value = 9750
$ws.Range("L69").Value = 29250
$ws.Range("N69").Value = -30872
$ws.Range("H70").Value = 400
$ws.Range("J70").Value = 400
$ws.Range("L70").Value = 1200
$ws.Range("N70").Value = -1830
$ws.Range("H71").Value = 85223.586
$ws.Range("J71").Value = 92743.91
$ws.Range("L71").Value = 834695.1900000001
$ws.Range("N71").Value = -842807.1900000001
$ws.Range("H72").Value = 9750
$ws.Range("J72").Value = 9750
$ws.Range("L72").Value = 87750
$ws.Range("N72").Value = -95862
$ws.Range("H73").Value = 400
$ws.Range("J73").Value = 400
$ws.Range("L73").Value = 1200
$ws.Range("N73").Value = -3384
$ws.Range("H92").Value = 555.8461
$ws.Range("I92").Value = 510.42856
$ws.Range("K92").Value = 1531.28568
$ws.Range("M92").Value = -283.28568
$ws.Range("H117").Value = 1070.4
$ws.Range("I117").Value = 1063
$ws.Range("J117").Value = 1100
$ws.Range("K117").Value = 3189
$ws.Range("L117").Value = 3300
$ws.Range("M117").Value = 253
$ws.Range("N117").Value = -10184

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 11195.272
$ws.Range("J93").Value = 11195.272
$ws.Range("L93").Value = 11195.272
$ws.Range("N93").Value = -14939.272
$ws.Range("H109").Value = 22551.455
$ws.Range("I109").Value = 7500
$ws.Range("K109").Value = 7500
$ws.Range("M109").Value = -6460
$ws.Range("H135").Value = 39976
$ws.Range("J135").Value = 39976
$ws.Range("L135").Value = 39976
$ws.Range("N135").Value = -50116
$ws.Range("H140").Value = 94552.37
$ws.Range("J140").Value = 94957.60000000001
$ws.Range("L140").Value = 94957.60000000001
$ws.Range("N140").Value = -105317.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4750
$ws.Range("I10").Value = 4500
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 4500
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -4360
$ws.Range("N10").Value = -5280
$ws.Range("H22").Value = 39375.25
$ws.Range("J22").Value = 75000
$ws.Range("L22").Value = 75000
$ws.Range("N22").Value = -75590
$ws.Range("H27").Value = 39375.25
$ws.Range("J27").Value = 75000
$ws.Range("L27").Value = 75000
$ws.Range("N27").Value = -75214
$ws.Range("H40").Value = 6947177
$ws.Range("J40").Value = 18521688
$ws.Range("L40").Value = 18521688
$ws.Range("N40").Value = -18521960
$ws.Range("H46").Value = 6375.3687
$ws.Range("I46").Value = 7795.533
$ws.Range("J46").Value = 1049.75
$ws.Range("K46").Value = 7795.533
$ws.Range("L46").Value = 1049.75
$ws.Range("M46").Value = -7607.533
$ws.Range("N46").Value = -1425.75
$ws.Range("H55").Value = 5061.4
$ws.Range("J55").Value = 6008.2
$ws.Range("L55").Value = 6008.2
$ws.Range("N55").Value = -6354.2
$ws.Range("H123").Value = 71170.625
$ws.Range("J123").Value = 74139.28999999999
$ws.Range("L123").Value = 74139.28999999999
$ws.Range("N123").Value = -83939.28999999999
$ws.Range("H129").Value = 70695
$ws.Range("J129").Value = 65000
$ws.Range("L129").Value = 65000
$ws.Range("N129").Value = -75000

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H136").Value = 1611.1818
$ws.Range("I136").Value = 1252.5555
$ws.Range("K136").Value = 3757.6665
$ws.Range("M136").Value = -1207.6665
